$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2
$ws.Range('D2').Value = '65.856.28'
$ws.Range('E2').Value = '  -1.37%  '
# Row 3
$ws.Range('D3').Value = '3.433.10'
$ws.Range('E3').Value = '  -1.05%  '
# Row 4
$c = $ws.Range('D4')
$c.NumberFormat = '@'
$c.Value = '0.999'
$c.ClearFormats()
$ws.Range('E4').Value = '  -0.19%  '
# Row 5
$c = $ws.Range('D5')
$c.NumberFormat = '@'
$c.Value = '597.13'
$c.ClearFormats()
$ws.Range('E5').Value = '  -0.26%  '
# Row 6
$c = $ws.Range('D6')
$c.NumberFormat = '@'
$c.Value = '142.15'
$c.ClearFormats()
$ws.Range('E6').Value = '  -3.05%  '
# Row 7
$ws.Range('D7').Value = '3.439.65'
$ws.Range('E7').Value = '  -0.77%  '
# Row 9
$c = $ws.Range('D9')
$c.NumberFormat = '@'
$c.Value = '0.471'
$c.ClearFormats()
$ws.Range('E9').Value = '  -1.24%  '
# Row 10
$c = $ws.Range('D10')
$c.NumberFormat = '@'
$c.Value = '8.06'
$c.ClearFormats()
$ws.Range('E10').Value = '  +6.72%  '
# Row 11
$c = $ws.Range('D11')
$c.NumberFormat = '@'
$c.Value = '0.134'
$c.ClearFormats()
$ws.Range('E11').Value = '  -4.92%  '
# Row 12
$c = $ws.Range('D12')
$c.NumberFormat = '@'
$c.Value = '0.407'
$c.ClearFormats()
$ws.Range('E12').Value = '  -3.21%  '
# Row 13
$ws.Range('D13').Value = '4.010.69'
$ws.Range('E13').Value = '  -1.13%  '
# Row 14
$c = $ws.Range('D14')
$c.NumberFormat = '@'
$c.Value = '0.0000200'
$c.ClearFormats()
$ws.Range('E14').Value = '  -5.21%  '
# Row 15
$c = $ws.Range('D15')
$c.NumberFormat = '@'
$c.Value = '29.73'
$c.ClearFormats()
$ws.Range('E15').Value = '  -4.51%  '
# Row 16
$ws.Range('B16').Value = 'WrappedEther'
$ws.Range('C16').Value = 'https://coinranking.com/coin/Mtfb0obXVh59u+wrappedether-weth'
$ws.Range('D16').Value = '3.423.63'
$ws.Range('E16').Value = '  -1.73%  '
# Row 17
$ws.Range('B17').Value = 'TRON'
$ws.Range('C17').Value = 'https://coinranking.com/coin/qUhEFk1I61atv+tron-trx'
$c = $ws.Range('D17')
$c.NumberFormat = '@'
$c.Value = '0.116'
$c.ClearFormats()
$ws.Range('E17').Value = '  -0.52%  '
# Row 18
$ws.Range('B18').Value = 'WrappedBTC'
$ws.Range('C18').Value = 'https://coinranking.com/coin/x4WXHge-vvFY+wrappedbtc-wbtc'
$ws.Range('D18').Value = '65.836.58'
$ws.Range('E18').Value = '  -1.36%  '
# Row 19
$c = $ws.Range('D19')
$c.NumberFormat = '@'
$c.Value = '10.31'
$c.ClearFormats()
$ws.Range('E19').Value = '  +3.28%  '
# Row 20
$c = $ws.Range('D20')
$c.NumberFormat = '@'
$c.Value = '6.10'
$c.ClearFormats()
$ws.Range('E20').Value = '  -4.20%  '
# Row 21
$c = $ws.Range('D21')
$c.NumberFormat = '@'
$c.Value = '14.59'
$c.ClearFormats()
$ws.Range('E21').Value = '  -4.15%  '
# Row 22
$c = $ws.Range('D22')
$c.NumberFormat = '@'
$c.Value = '415.80'
$c.ClearFormats()
$ws.Range('E22').Value = '  -3.71%  '
# Row 23
$c = $ws.Range('D23')
$c.NumberFormat = '@'
$c.Value = '0.578'
$c.ClearFormats()
$ws.Range('E23').Value = '  -4.08%  '
# Row 24
$c = $ws.Range('D24')
$c.NumberFormat = '@'
$c.Value = '77.19'
$c.ClearFormats()
$ws.Range('E24').Value = '  -2.41%  '
# Row 25
$ws.Range('E25').Value = '  +0.07%  '
# Row 26
$c = $ws.Range('D26')
$c.NumberFormat = '@'
$c.Value = '0.0000111'
$c.ClearFormats()
$ws.Range('E26').Value = '  -7.16%  '
# Row 27
$c = $ws.Range('D27')
$c.NumberFormat = '@'
$c.Value = '9.30'
$c.ClearFormats()
$ws.Range('E27').Value = '  -4.29%  '
# Row 28
$c = $ws.Range('D28')
$c.NumberFormat = '@'
$c.Value = '7.95'
$c.ClearFormats()
$ws.Range('E28').Value = '  -4.60%  '
# Row 29
$c = $ws.Range('D29')
$c.NumberFormat = '@'
$c.Value = '2.43'
$c.ClearFormats()
$ws.Range('E29').Value = '  -1.43%  '
# Row 30
$ws.Range('E30').Value = '  -0.08%  '
# Row 31
$c = $ws.Range('D31')
$c.NumberFormat = '@'
$c.Value = '0.161'
$c.ClearFormats()
$ws.Range('E31').Value = '  -3.58%  '
# Row 32
$c = $ws.Range('D32')
$c.NumberFormat = '@'
$c.Value = '1.46'
$c.ClearFormats()
$ws.Range('E32').Value = '  -7.22%  '
# Row 33
$c = $ws.Range('D33')
$c.NumberFormat = '@'
$c.Value = '24.65'
$c.ClearFormats()
$ws.Range('E33').Value = '  -2.09%  '
# Row 34
$ws.Range('D34').Value = '3.427.58'
$ws.Range('E34').Value = '  -0.95%  '
# Row 35
$ws.Range('E35').Value = '  -0.07%  '
# Row 36
$c = $ws.Range('D36')
$c.NumberFormat = '@'
$c.Value = '1.69'
$c.ClearFormats()
$ws.Range('E36').Value = '  -5.21%  '
# Row 37
$c = $ws.Range('D37')
$c.NumberFormat = '@'
$c.Value = '5.51'
$c.ClearFormats()
$ws.Range('E37').Value = '  -6.16%  '
# Row 38
$c = $ws.Range('D38')
$c.NumberFormat = '@'
$c.Value = '7.53'
$c.ClearFormats()
$ws.Range('E38').Value = '  -3.70%  '
# Row 39
$c = $ws.Range('D39')
$c.NumberFormat = '@'
$c.Value = '0.997'
$c.ClearFormats()
$ws.Range('E39').Value = '  -0.19%  '
# Row 40
$c = $ws.Range('D40')
$c.NumberFormat = '@'
$c.Value = '169.13'
$c.ClearFormats()
$ws.Range('E40').Value = '  -3.21%  '
# Row 41
$c = $ws.Range('D41')
$c.NumberFormat = '@'
$c.Value = '0.0858'
$c.ClearFormats()
$ws.Range('E41').Value = '  -1.90%  '
# Row 42
$c = $ws.Range('D42')
$c.NumberFormat = '@'
$c.Value = '0.874'
$c.ClearFormats()
$ws.Range('E42').Value = '  -1.79%  '
# Row 43
$c = $ws.Range('D43')
$c.NumberFormat = '@'
$c.Value = '5.06'
$c.ClearFormats()
$ws.Range('E43').Value = '  -5.90%  '
# Row 44
$c = $ws.Range('D44')
$c.NumberFormat = '@'
$c.Value = '1.90'
$c.ClearFormats()
$ws.Range('E44').Value = '  -9.54%  '
# Row 45
$c = $ws.Range('D45')
$c.NumberFormat = '@'
$c.Value = '45.42'
$c.ClearFormats()
$ws.Range('E45').Value = '  -1.71%  '
# Row 46
$c = $ws.Range('D46')
$c.NumberFormat = '@'
$c.Value = '26.25'
$c.ClearFormats()
$ws.Range('E46').Value = '  -7.91%  '
# Row 47
$ws.Range('E47').Value = '  -2.08%  '
# Row 48
$c = $ws.Range('D48')
$c.NumberFormat = '@'
$c.Value = '7.06'
$c.ClearFormats()
$ws.Range('E48').Value = '  -4.44%  '
# Row 49
$c = $ws.Range('D49')
$c.NumberFormat = '@'
$c.Value = '2.29'
$c.ClearFormats()
$ws.Range('E49').Value = '  -4.78%  '
# Row 50
$c = $ws.Range('D50')
$c.NumberFormat = '@'
$c.Value = '0.922'
$c.ClearFormats()
$ws.Range('E50').Value = '  -4.89%  '
# Row 51
$c = $ws.Range('D51')
$c.NumberFormat = '@'
$c.Value = '0.232'
$c.ClearFormats()
$ws.Range('E51').Value = '  -4.52%  '
